# Update Wnt1-Ror2 LR-pairs sheet with new TPM-derived values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last three data rows (old rows 5-7); the table shrinks from
# 6 data rows to 3 data rows (A1:T7 -> A1:T4).
$ws.Range("A5:T7").EntireRow.Delete() | Out-Null

# Row 2: FAPs | Wnt1 | Ror2 | ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Wnt1"
$ws.Range("C2").Value = "Ror2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.03910466666666667
$ws.Range("H2").Value = 0.117314
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.491263
$ws.Range("N2").Value = 1.473789
$ws.Range("O2").Value = 0.08507192955174298
$ws.Range("P2").Value = 0.08507192955174299
$ws.Range("Q2").Value = 0.01921067586066667
$ws.Range("R2").Value = 0.172896082746
$ws.Range("S2").Value = 0.08507192955174298
$ws.Range("T2").Value = 0.08507192955174299

# Row 3: FAPs | Wnt1 | Ror2 | FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Wnt1"
$ws.Range("C3").Value = "Ror2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.03910466666666667
$ws.Range("H3").Value = 0.117314
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.669265333333333
$ws.Range("N3").Value = 14.007796
$ws.Range("O3").Value = 0.8085758778815603
$ws.Range("P3").Value = 0.8085758778815605
$ws.Range("Q3").Value = 0.1825900644382222
$ws.Range("R3").Value = 1.643310579944
$ws.Range("S3").Value = 0.8085758778815603
$ws.Range("T3").Value = 0.8085758778815605

# Row 4: FAPs | Wnt1 | Ror2 | MuSCs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Wnt1"
$ws.Range("C4").Value = "Ror2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.03910466666666667
$ws.Range("H4").Value = 0.117314
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.6141496666666667
$ws.Range("N4").Value = 1.842449
$ws.Range("O4").Value = 0.1063521925666967
$ws.Range("P4").Value = 0.1063521925666967
$ws.Range("Q4").Value = 0.02401611799844445
$ws.Range("R4").Value = 0.216145061986
$ws.Range("S4").Value = 0.1063521925666967
$ws.Range("T4").Value = 0.1063521925666967
